# Append six more "journal" rows (28-33), continuing the existing pattern
# found in rows 4-27: column A holds a sequential number, columns B-G cycle
# through two alternating templates of text values. Row 27's formatting is
# copied down to each new row so the shared cell style (s="9") is reused
# instead of minting new style entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateOdd  = @("Тимків Віталій Дмитрович", "А0000", "Тимків Дмитро Віталійович", "Уганда, гасити вагнерів", "01.01.2025 Краківець", "01.01.2026 Подобовець")
$templateEven = @("Тимків Віталій Дмитрович 1", "А0000 1", "Тимків Дмитро Віталійович 1", "Уганда, гасити вагнерів 1", "01.01.2025 Краківець 1", "01.01.2026 Подобовець 1")

$startRow = 28
$endRow = 33
$seq = 25

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Range("A" + ($r - 1) + ":G" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $seq

    if ((($r - $startRow) % 2) -eq 0) {
        $template = $templateOdd
    } else {
        $template = $templateEven
    }

    for ($i = 0; $i -lt $template.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $template[$i]
    }

    $seq++
}

$excel.CutCopyMode = 0
